$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("positions")

# --- Education row (row 2): change separator between degrees from "|" to "&" ---
$ws.Range("C2").Value = "M.S., Decision Analytics        &       Master of Business Administration"

# --- Markel Corporation / Associate Data Scientist row (row 6): rewrite description bullets ---
$ws.Range("H6").Value = "Data scientist on an agile sprint team that builds, deploys, and maintains production applications that help prioritize and monitor incoming business"
$ws.Range("I6").Value = "Develops, maintains, and deploys machine learning models to end users and conducts ad hoc analysis that help drive decision making"

# --- Move the active selection on the "positions" sheet to K7 ---
$ws.Range("K7").Select()
